# Add daily power updates: append rows 39-41 (2018-09-20, 2018-09-21,
# 2018-09-22) to the comforter-cda sheet's data table, growing the table
# and dimension/selection to match, mirroring the formulas used in the
# existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comforter-cda")

# Row 39 (2018-09-20): no start/end time recorded yet; duration formulas
# are present but evaluate against blank B/C so they read as 0.
$ws.Range("A39").Value = 43363
$ws.Range("D39").Formula = "=(C39-B39)* 1440"
$ws.Range("E39").Formula = "=IF(C39>B39, (C39-B39)*1440, (B39-C39)*1440)"
$ws.Range("F39").Formula = "=ABS((C39-B39)*1440)"

# Row 40 (2018-09-21): same situation - date only, no times yet.
$ws.Range("A40").Value = 43364
$ws.Range("D40").Formula = "=(C40-B40)* 1440"
$ws.Range("E40").Formula = "=IF(C40>B40, (C40-B40)*1440, (B40-C40)*1440)"
$ws.Range("F40").Formula = "=ABS((C40-B40)*1440)"

# Row 41 (2018-09-22): full entry with start and end time.
$ws.Range("A41").Value = 43365
$ws.Range("B41").Value = 0.79583333333333339
$ws.Range("C41").Value = 0.99930555555555556
$ws.Range("D41").Formula = "=(C41-B41)* 1440"
$ws.Range("E41").Formula = "=IF(C41>B41, (C41-B41)*1440, (B41-C41)*1440)"
$ws.Range("F41").Formula = "=ABS((C41-B41)*1440)"

# Grow the table over the new rows - this also carries each column's
# formatting (date / time / duration styles) onto the new rows.
$ws.ListObjects.Item("comforter_cda_table").Resize($ws.Range("A1:F41"))

# Match the updated view state: selection on the new last cell and the
# sheet's used-range dimension.
$ws.Range("D41").Select()
